$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the date serial value in I1 (33121 -> 33133)
$ws.Range("I1").Value = 33133

# Clear the "x" marker + bill amount for ComCast Business (row 3)
$ws.Range("B3").ClearContents()
$ws.Range("D3").ClearContents()

# Clear the "x" marker + bill amount for Hinckley Springs (row 8)
$ws.Range("B8").ClearContents()
$ws.Range("D8").ClearContents()

# Clear the "x" marker + bill amount for Pitney Bowes Purchase Power (row 12)
$ws.Range("B12").ClearContents()
$ws.Range("D12").ClearContents()

# Mark Salvi Salvi & Wifler (row 13) as "x" with a new bill amount
$ws.Range("B13").Value = "x"
$ws.Range("D13").Value = 343.65

# Update the active selection to match the author's final cursor position
$ws.Range("C13").Select()
